$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first 16 data rows (rows 2-17, the oldest observations) so the
# series starts at the next release date; remaining rows shift up and the
# newest rows (now unused at the bottom) disappear along with the old dimension.
$ws.Range("A2:B17").EntireRow.Delete()
